$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row above current row 8 so the existing rows 8-15 shift
# down to 9-16, matching the diff's new dimension A1:F16.
$ws.Rows.Item(8).Insert()

# Copy the date/time number format used by column F (from row 7, which is
# untouched) onto the newly inserted row's F cell so style "2" is preserved.
$ws.Cells.Item(7, 6).Copy()
$ws.Cells.Item(8, 6).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-write the full block of rows 8-16 with the corrected / re-ordered data.
$data = @(
    @(761, "Mate Security", "VP Sales - US",         "Michael Persechini", "1st Interview", 45995),
    @(769, "TraceBit",      "VP Sales NYC",            "Matthew Schaner",    "CV Sent",       45987),
    @(770, "TraceBit",      "Account Executive NYC",   "Ben Brighton",       "4th Interview",  45992),
    @(770, "TraceBit",      "Account Executive NYC",   "Bruna Corbin",       "1st Interview",  45986),
    @(770, "TraceBit",      "Account Executive NYC",   "Jameson Schwartz",   "1st Interview",  45982),
    @(770, "TraceBit",      "Account Executive NYC",   "Ryan Drillock",      "2nd Interview",  45987),
    @(791, "Adaptive6",     "Head of Sales (US)",      "Bryan Pierrot",      "CV Sent",       45987),
    @(791, "Adaptive6",     "Head of Sales (US)",      "Dan Baldassano",     "CV Sent",       45987),
    @(834, "Blockaid",      "Regional Director US",    "Jeff White",         "1st Interview",  45999)
)

$r = 8
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
